$d = $word.ActiveDocument

$d.Content.Find.Execute("265×6=1590", $true, $false, $false, $false, $false, $true, 1, $false, "195×6=1170", 2) | Out-Null
$d.Content.Find.Execute("426×7=2982", $true, $false, $false, $false, $false, $true, 1, $false, "139×9=1251", 2) | Out-Null
$d.Content.Find.Execute("915×5=4575", $true, $false, $false, $false, $false, $true, 1, $false, "493×4=1972", 2) | Out-Null
$d.Content.Find.Execute("107×2=214", $true, $false, $false, $false, $false, $true, 1, $false, "293×4=1172", 2) | Out-Null
$d.Content.Find.Execute("273×2=546", $true, $false, $false, $false, $false, $true, 1, $false, "407×9=3663", 2) | Out-Null
$d.Content.Find.Execute("173×9=1557", $true, $false, $false, $false, $false, $true, 1, $false, "458×4=1832", 2) | Out-Null
$d.Content.Find.Execute("605×5=3025", $true, $false, $false, $false, $false, $true, 1, $false, "450×6=2700", 2) | Out-Null
$d.Content.Find.Execute("274×4=1096", $true, $false, $false, $false, $false, $true, 1, $false, "355×2=710", 2) | Out-Null
$d.Content.Find.Execute("708×2=1416", $true, $false, $false, $false, $false, $true, 1, $false, "869×8=6952", 2) | Out-Null
$d.Content.Find.Execute("426×2=852", $true, $false, $false, $false, $false, $true, 1, $false, "748×8=5984", 2) | Out-Null
$d.Content.Find.Execute("232×6=1392", $true, $false, $false, $false, $false, $true, 1, $false, "621×6=3726", 2) | Out-Null
$d.Content.Find.Execute("526×4=2104", $true, $false, $false, $false, $false, $true, 1, $false, "298×6=1788", 2) | Out-Null
$d.Content.Find.Execute("548×5=2740", $true, $false, $false, $false, $false, $true, 1, $false, "416×7=2912", 2) | Out-Null
$d.Content.Find.Execute("778×2=1556", $true, $false, $false, $false, $false, $true, 1, $false, "725×7=5075", 2) | Out-Null
$d.Content.Find.Execute("559×7=3913", $true, $false, $false, $false, $false, $true, 1, $false, "684×8=5472", 2) | Out-Null
$d.Content.Find.Execute("803×7=5621", $true, $false, $false, $false, $false, $true, 1, $false, "172×8=1376", 2) | Out-Null
$d.Content.Find.Execute("851×4=3404", $true, $false, $false, $false, $false, $true, 1, $false, "135×3=405", 2) | Out-Null
$d.Content.Find.Execute("123×6=738", $true, $false, $false, $false, $false, $true, 1, $false, "740×8=5920", 2) | Out-Null
$d.Content.Find.Execute("457×7=3199", $true, $false, $false, $false, $false, $true, 1, $false, "137×9=1233", 2) | Out-Null
$d.Content.Find.Execute("419×3=1257", $true, $false, $false, $false, $false, $true, 1, $false, "670×3=2010", 2) | Out-Null
$d.Content.Find.Execute("809×9=7281", $true, $false, $false, $false, $false, $true, 1, $false, "981×3=2943", 2) | Out-Null
$d.Content.Find.Execute("887×4=3548", $true, $false, $false, $false, $false, $true, 1, $false, "143×3=429", 2) | Out-Null
$d.Content.Find.Execute("895×4=3580", $true, $false, $false, $false, $false, $true, 1, $false, "619×6=3714", 2) | Out-Null
$d.Content.Find.Execute("961×9=8649", $true, $false, $false, $false, $false, $true, 1, $false, "830×7=5810", 2) | Out-Null
$d.Content.Find.Execute("675×6=4050", $true, $false, $false, $false, $false, $true, 1, $false, "442×5=2210", 2) | Out-Null
